$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns B:E to plain Text format so that numeric-looking strings
# (e.g. "101.01", "1.00") are preserved as text rather than being
# auto-converted to numbers, matching the source inlineStr cells.
$ws.Range("B2:E51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "50.881.14"
$ws.Range("E2").Value = "  -0.87%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.926.82"
$ws.Range("E3").Value = "  -1.47%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.28%  "

# Row 5 - BNB
$ws.Range("D5").Value = "373.13"
$ws.Range("E5").Value = "  -1.82%  "

# Row 6 - Solana
$ws.Range("D6").Value = "101.01"
$ws.Range("E6").Value = "  -4.14%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -1.30%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.15%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.579"
$ws.Range("E9").Value = "  -2.35%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "35.79"
$ws.Range("E10").Value = "  -4.03%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -0.68%  "

# Row 12 - Dogecoin
$ws.Range("E12").Value = "  -0.43%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.391.71"
$ws.Range("E13").Value = "  -1.26%  "

# Row 14 - Uniswap (was Chainlink)
$ws.Range("B14").Value = "Uniswap"
$ws.Range("C14").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D14").Value = "12.28"
$ws.Range("E14").Value = "  +66.34%  "

# Row 15 - Chainlink (was Polkadot)
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "17.85"
$ws.Range("E15").Value = "  -3.03%  "

# Row 16 - Polkadot (was Uniswap)
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "7.34"
$ws.Range("E16").Value = "  -2.80%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.941.14"
$ws.Range("E17").Value = "  -0.43%  "

# Row 18 - Polygon
$ws.Range("D18").Value = "0.964"
$ws.Range("E18").Value = "  -0.05%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "50.848.92"
$ws.Range("E19").Value = "  -0.88%  "

# Row 20 - ImmutableX
$ws.Range("E20").Value = "  -6.38%  "

# Row 21 - InternetComputer(DFINITY)
$ws.Range("D21").Value = "12.35"
$ws.Range("E21").Value = "  -4.34%  "

# Row 22 - ShibaInu
$ws.Range("D22").Value = "0.0₃0949"
$ws.Range("E22").Value = "  -1.39%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "263.31"
$ws.Range("E23").Value = "  +0.76%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "68.28"
$ws.Range("E24").Value = "  -1.61%  "

# Row 25 - PancakeSwap
$ws.Range("D25").Value = "3.16"
$ws.Range("E25").Value = "  +12.11%  "

# Row 26 - Filecoin
$ws.Range("E26").Value = "  +4.02%  "

# Row 27 - RenderToken
$ws.Range("D27").Value = "7.47"
$ws.Range("E27").Value = "  -0.24%  "

# Row 28 - Dai
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.02%  "

# Row 29 - Kaspa
$ws.Range("D29").Value = "0.165"
$ws.Range("E29").Value = "  -3.75%  "

# Row 30 - EthereumClassic
$ws.Range("D30").Value = "25.40"
$ws.Range("E30").Value = "  -1.83%  "

# Row 31 - Hedera
$ws.Range("D31").Value = "0.109"
$ws.Range("E31").Value = "  -2.83%  "

# Row 32 - Cosmos
$ws.Range("D32").Value = "9.93"
$ws.Range("E32").Value = "  +0.36%  "

# Row 33 - OKB
$ws.Range("D33").Value = "50.54"
$ws.Range("E33").Value = "  -1.31%  "

# Row 34 - Toncoin
$ws.Range("E34").Value = "  -2.95%  "

# Row 35 - VeChain
$ws.Range("D35").Value = "0.0437"
$ws.Range("E35").Value = "  -1.97%  "

# Row 36 - InjectiveProtocol
$ws.Range("D36").Value = "32.63"
$ws.Range("E36").Value = "  -7.36%  "

# Row 37 - FirstDigitalUSD
$ws.Range("E37").Value = "  +0.07%  "

# Row 38 - LidoDAOToken
$ws.Range("D38").Value = "3.11"
$ws.Range("E38").Value = "  +1.78%  "

# Row 39 - Stellar
$ws.Range("E39").Value = "  -0.91%  "

# Row 40 - Celestia
$ws.Range("D40").Value = "16.07"
$ws.Range("E40").Value = "  -6.68%  "

# Row 41 - Stacks
$ws.Range("D41").Value = "2.46"
$ws.Range("E41").Value = "  -5.23%  "

# Row 42 - ARBITRUM
$ws.Range("E42").Value = "  -4.36%  "

# Row 43 - Monero
$ws.Range("E43").Value = "  -3.43%  "

# Row 44 - EnergySwap
$ws.Range("D44").Value = "20.91"
$ws.Range("E44").Value = "  -4.88%  "

# Row 45 - WEMIXToken (was TheGraph)
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "2.04"
$ws.Range("E45").Value = "  -1.03%  "

# Row 46 - TheGraph (was WEMIXToken)
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Value = "0.273"
$ws.Range("E46").Value = "  -5.86%  "

# Row 47 - NEARProtocol
$ws.Range("D47").Value = "3.26"
$ws.Range("E47").Value = "  +0.63%  "

# Row 48 - ApeXProtocol
$ws.Range("E48").Value = "  -3.11%  "

# Row 49 - Maker
$ws.Range("D49").Value = "1.985.97"
$ws.Range("E49").Value = "  -2.78%  "

# Row 50 - BEAM
$ws.Range("D50").Value = "0.0327"
$ws.Range("E50").Value = "  -4.53%  "

# Row 51 - TrustWalletToken
$ws.Range("E51").Value = "  -0.07%  "

# Restore the default (unstyled) cell style on the data cells we touched,
# since the source workbook keeps these cells on the default style (no 's'
# attribute) even though their content is textual.
$ws.Range("B2:E51").Style = "Normal"

